$d = $word.ActiveDocument

# Locate the reference paragraph that must be kept ("SERAFINI, ...").
# Everything between it and the page-break paragraph near the end of the
# document (an empty paragraph, the "Ver no Jupiter ..." line, and the
# "© 2020 ..." footer line) should be removed, while the trailing empty
# paragraph is left intact.
$anchor = $d.Content
$anchor.Find.Execute("SERAFINI, Maria José. Como escrever textos. 5.ed. São Paulo: Globo, 1992.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$foundEnd = $anchor.End

# Resolve which document paragraph (1-based Paragraphs index) contains the
# match just found, walking the real paragraph collection so indices line
# up with true paragraph-mark boundaries.
$paraCount = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le ($foundEnd - 1) -and $p.Range.End -ge $foundEnd) {
        $anchorIndex = $i
        break
    }
}

# The three paragraphs to delete are the ones immediately following the
# anchor paragraph: the blank line, "Ver no Jupiter ...", and "© 2020 ...".
$firstToDelete = $d.Paragraphs.Item($anchorIndex + 1)
$lastToDelete = $d.Paragraphs.Item($anchorIndex + 3)

$deleteRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
$deleteRange.Delete()
